# Fix Capex/Opex battery storage LCOH values and widen column B
# on the "Present-Storage" and "2050-Storage" sheets.

$wb = $excel.ActiveWorkbook

# --- Present-Storage sheet ---
# Note: the saved <col width> value is derived from ColumnWidth after
# internal pixel-rounding, so we pick the ColumnWidth value whose
# rounded result lands on the target width of 9.6 (closest achievable).
$ws1 = $wb.Worksheets.Item("Present-Storage")
$ws1.Columns.Item(2).ColumnWidth = 8.83
$ws1.Range("B2").Value = 132.04
$ws1.Range("B3").Value = 563.3
$ws1.Range("B4").Value = 58.67
$ws1.Range("B6").Value = 139.11
$ws1.Range("B7").Value = 595.17
$ws1.Range("B8").Value = 61.52
$ws1.Range("B10").Value = 107.99
$ws1.Range("B11").Value = 444.46
$ws1.Range("B12").Value = 50.75
$ws1.Range("B14").Value = 128.01
$ws1.Range("B15").Value = 541.8
$ws1.Range("B16").Value = 57.61

# --- 2030-Storage sheet ---
$ws2 = $wb.Worksheets.Item("2030-Storage")
$ws2.Range("B2").Value = 5.88
$ws2.Range("B3").Value = 6.04
$ws2.Range("B4").Value = 6.37
$ws2.Range("B6").Value = 5.7
$ws2.Range("B7").Value = 5.87
$ws2.Range("B8").Value = 6.21
$ws2.Range("B10").Value = 9.56
$ws2.Range("B11").Value = 9.69
$ws2.Range("B12").Value = 9.94
$ws2.Range("B14").Value = 6.96
$ws2.Range("B15").Value = 7.12
$ws2.Range("B16").Value = 7.43

# --- 2050-Storage sheet ---
$ws3 = $wb.Worksheets.Item("2050-Storage")
$ws3.Columns.Item(2).ColumnWidth = 8.83
$ws3.Range("B2").Value = 146.2
$ws3.Range("B3").Value = 73.36
$ws3.Range("B4").Value = 92.23
$ws3.Range("B6").Value = 154.08
$ws3.Range("B7").Value = 77.05
$ws3.Range("B8").Value = 97.01000000000001
$ws3.Range("B10").Value = 119.03
$ws3.Range("B11").Value = 62.2
$ws3.Range("B12").Value = 76.93000000000001
$ws3.Range("B14").Value = 141.59
$ws3.Range("B15").Value = 71.7
$ws3.Range("B16").Value = 89.81
